$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently marks the very last edit point
#    in the document (right after the "]" that closes the sample
#    "$inject" array, near the end of the file). Word always keeps
#    only one "_GoBack" bookmark, relocating it to wherever the most
#    recent edit happened. Remove the stale one now; we re-create it
#    below at the new edit location.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate "blueprint/app/NovaWeb." - the text that is immediately
#    followed by the screenshot picture in section 2. We need a
#    manual line break right after that "." and before the picture,
#    so the picture starts on its own line (the alignment fix).
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("blueprint/app/NovaWeb.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text 'blueprint/app/NovaWeb.'"
}
# $rng now spans just the matched text; collapse to its end, i.e. the
# boundary right after "." and right before the inline picture.
$rng.Collapse(0)
$breakPos = $rng.Start

# Word's COM bridge here mis-places a break inserted exactly on a
# boundary that touches an inline picture (it lands the break after
# the picture instead of before it). Work around this by temporarily
# pushing the picture away with placeholder text, inserting the break
# where it's now safely surrounded by plain text, then removing the
# placeholder again.
$placeholder = "ZZPLACEHOLDERZZ"
$rng.InsertAfter($placeholder)

$breakRange = $d.Range($breakPos, $breakPos)
$breakRange.InsertBreak(6)

$placeholderRange = $d.Range($breakPos + 1, $breakPos + 1 + $placeholder.Length)
$placeholderRange.Delete()

# ------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark, collapsed, right after the new
#    line break and right before the picture.
# ------------------------------------------------------------------
$bookmarkPos = $breakPos + 1
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
